# Daily attendance processing - reverse the order of names listed in the
# "Recorded By" column (column G) wherever multiple recorders are present
# (i.e. the value contains a comma-separated list of names/emails).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val.Split(",")
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }
        # Reverse the order of the names while keeping original text/casing
        $reversed = $trimmed[($trimmed.Length - 1)..0]
        $newVal = [string]::Join(", ", $reversed)
        $cell.Value2 = $newVal
    }
}
